$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Break the merges that must change shape before we touch their contents.
# ---------------------------------------------------------------------------
$ws.Range("A2:A20").UnMerge()
$ws.Range("C2:E2").UnMerge()
$ws.Range("C35:D35").UnMerge()
$ws.Range("C36:D39").UnMerge()

# ---------------------------------------------------------------------------
# 2. Row 2 / Row 3 pinout header table.
# ---------------------------------------------------------------------------
# Row 2 group headers
$ws.Range("C2").Value = "BNO055IMU"
$ws.Range("I2").Value = "Button"
$ws.Range("K2").Value = "SD Adapter"
$ws.Range("Q2").Value = "NRF24l01 standard"

# Row 3 pin headers
$ws.Range("C3").Value = "Vin"
$ws.Range("D3").Value = "3vo"
$ws.Range("E3").Value = "GND"
$ws.Range("F3").Value = "SDA"
$ws.Range("G3").Value = "SCL"
$ws.Range("H3").Value = "RST"
$ws.Range("I3").Value = "side1"
$ws.Range("J3").Value = "side2"
$ws.Range("K3").Value = "CS"
$ws.Range("L3").Value = "SCK"
$ws.Range("M3").Value = "MOSI"
$ws.Range("N3").Value = "MISO"
$ws.Range("O3").Value = "VCC"
$ws.Range("P3").Value = "GND"

# ---------------------------------------------------------------------------
# 3. Key / legend block (rows 35-39) text shuffled down by one row.
# ---------------------------------------------------------------------------
$ws.Range("C35").ClearContents()
$ws.Range("C36").Value = "Key"
$ws.Range("C37").Value = "Microcontroller pin numbers are according to "

# ---------------------------------------------------------------------------
# 4. Re-merge the ranges into their new shapes.
# ---------------------------------------------------------------------------
$ws.Range("A2:A6").Merge()
$ws.Range("C2:H2").Merge()
$ws.Range("I2:J2").Merge()
$ws.Range("K2:P2").Merge()
$ws.Range("Q2:W2").Merge()
$ws.Range("C36:D36").Merge()
$ws.Range("C37:D39").Merge()

# ---------------------------------------------------------------------------
# 5. Formatting.
# ---------------------------------------------------------------------------
# A2:A6 "Microcontroller" label: centered + rotated 90 degrees.
$r = $ws.Range("A2:A6")
$r.HorizontalAlignment = -4108
$r.VerticalAlignment = -4108
$r.Orientation = 90

# A7:A20 filler cells: drop back to plain vertical-centering (no longer part
# of the big merge / no longer centered horizontally).
$ws.Range("A7:A20").HorizontalAlignment = 1

# New header cells on rows 2 & 3 need the same centered look as their peers.
$ws.Range("C2:W2").HorizontalAlignment = -4108
$ws.Range("C2:W2").VerticalAlignment = -4108
$ws.Range("C3:P3").HorizontalAlignment = -4108
$ws.Range("C3:P3").VerticalAlignment = -4108

# C35:D35 filler cells: plain vertical-centering only.
$ws.Range("C35:D35").HorizontalAlignment = 1

# C36:D36 "Key" label: centered, no wrap/rotation.
$r2 = $ws.Range("C36:D36")
$r2.HorizontalAlignment = -4108
$r2.VerticalAlignment = -4108
$r2.WrapText = $false
$r2.Orientation = 0

# C37:D39 long description: centered + rotated 90 + wrapped, taller rows.
$r3 = $ws.Range("C37:D39")
$r3.HorizontalAlignment = -4108
$r3.VerticalAlignment = -4108
$r3.WrapText = $true
$r3.Orientation = 90

$ws.Rows.Item(37).RowHeight = 25.5
$ws.Rows.Item(38).RowHeight = 25.5
$ws.Rows.Item(39).RowHeight = 25.5

# ---------------------------------------------------------------------------
# 6. Selection cursor, matching the saved workbook view.
# ---------------------------------------------------------------------------
$ws.Range("S3").Select()
